# Apply the permutation of species-observation rows 4-8.
# Only columns A, B, D, E, F, G, H, Q, R differ between the affected rows;
# all other columns (C, I, K, P, S..AY) stay identical per row.
#
# The row data is cyclically permuted:
#   row 4 <-> row 6            (2-cycle)
#   row 5 -> row 7 -> row 8 -> row 5   (3-cycle)
#
# Capture all source values first, then write them out, so that writes
# to one row never clobber data still needed for another row.

$ws = $excel.ActiveWorkbook.ActiveSheet

$cols = @("A", "B", "D", "E", "F", "G", "H", "Q", "R")
$rows = 4..8

# Snapshot current values for each affected cell.
# (Use Value2 for the read - it returns the raw value; Value is used below
# for the write, per the documented Range.Value setter usage.)
$snapshot = @{}
foreach ($r in $rows) {
    foreach ($col in $cols) {
        $snapshot["$col$r"] = $ws.Range("$col$r").Value2
    }
}

# current row -> target row (where that row's data ends up)
$mapping = @{ 4 = 6; 5 = 7; 6 = 4; 7 = 8; 8 = 5 }

foreach ($srcRow in $rows) {
    $dstRow = $mapping[$srcRow]
    foreach ($col in $cols) {
        $ws.Range("$col$dstRow").Value = $snapshot["$col$srcRow"]
    }
}
